$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.929.26'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.94%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.815.63'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.09%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.82%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9990'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.33%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4695'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.46%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3704'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.16%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07369'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.17%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8714'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.55%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.42'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.64%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.861.89'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.96%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.362'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.35%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.511'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.04%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07064'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.52%  '

# Row 16
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.81'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.26%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.08%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008706'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.37%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9991'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.21%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.08%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.005.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.66%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.350'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.71%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.55'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.22%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.114.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.42%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.898'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.76%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.04%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.42'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.61%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.165'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.96%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.295'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.52%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.48'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.16%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08937'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.55%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7602'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.65%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.159'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.97%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.474'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.80%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.926'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.07%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9976'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.44%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.094'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.62%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01954'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.26%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05262'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.42%  '

# Row 40
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.934'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.84%  '

# Row 41
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5347'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.40%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.214'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.42%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.371'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.24%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1661'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.21%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.471'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.16%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4957'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.36%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.43%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.678'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.66%  '

# Row 49
$ws.Range('E49').Value = '  -0.41%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.26'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.76%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06285'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.58%  '

